# Add a new, hidden worksheet ("HiddenSheet") after the existing Sheet1,
# containing a single cell (A4) explaining that the sheet is hidden.
# This mirrors the upstream fix for "Added support for hidden pages" (#9)
# and the NPE fix (#8).

$wb = $excel.ActiveWorkbook

# Remember the currently active sheet so we can restore the original
# selection/active-sheet state once we're done.
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 (Worksheets.Add with no "Before"
# and an "After" sheet appends immediately following it).
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "HiddenSheet"

# Populate the explanatory text in A4.
$ws.Range("A4").Value = "This sheet is hidden, and shouldn't appear in the output."
[void]$ws.Range("A4").Select()

# Hide the sheet (regular hidden, not "very hidden").
$ws.Visible = $false

# Restore the original active sheet/selection.
[void]$sheet1.Select()
